# WcDonaldas_3Sprintas.pptx - Slide 7 "Gito repo" textbox update
#
# - Turns the existing GitHub URL run into a hyperlink.
# - Adds a blank line.
# - Adds a "Subuildintas projektas:" line (built from four separate runs,
#   matching how PowerPoint splits runs after an autocorrect/spellcheck pass).
# - Adds the Google Drive folder link as a new line.
# - Resizes the auto-fit textbox to its new (taller) height.

$p = $ppt.ActivePresentation

# Locate the slide / shape robustly by scanning for the known marker text,
# instead of hard-coding slide/shape indices.
$targetSlide = $null
$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $sh = $sl.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -like "*github.com/MartynasKul/WcDSim*") {
                $targetSlide = $sl
                $targetShape = $sh
            }
        }
    }
}

$sh = $targetShape
$cr = [char]13

# 1) Append the new paragraphs first (so formatting/hyperlink work below only
#    touches the text that should carry it).

# Blank paragraph.
$t = $sh.TextFrame.TextRange
$t.InsertAfter($cr) | Out-Null

# "Subuildintas projektas:" paragraph, written as four separate runs.
$t = $sh.TextFrame.TextRange
$t.InsertAfter($cr) | Out-Null
$t = $sh.TextFrame.TextRange
$t.InsertAfter("Subuildintas") | Out-Null
$t = $sh.TextFrame.TextRange
$t.InsertAfter(" ") | Out-Null
$t = $sh.TextFrame.TextRange
$t.InsertAfter("projektas") | Out-Null
$t = $sh.TextFrame.TextRange
$t.InsertAfter(":") | Out-Null

# Google Drive folder link paragraph.
$t = $sh.TextFrame.TextRange
$t.InsertAfter($cr) | Out-Null
$t = $sh.TextFrame.TextRange
$t.InsertAfter("https://drive.google.com/drive/folders/11h6vm97gGb5ItiLMrskLFJmY8HVqxrKb?usp=sharing") | Out-Null

# 2) Turn the first line's text into a hyperlink pointing at itself.
$tr = $sh.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.ActionSettings(1).Hyperlink.Address = "https://github.com/MartynasKul/WcDSim" | Out-Null

# 3) The textbox auto-fits to its text; grow it to match the new line count.
$sh.Height = 92.0907
